$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{}
$values[2] = 'Amen amen amen'
$values[3] = 'I don''t know who you are or what your name is but if you''re Reading this, I send healing to your entire body in the name of Jesus'
$values[4] = 'Amén '
$values[5] = 'Breaking News God is going to Double your Blessings this week, get ready to Receive them.'
$values[6] = 'Casper Miguel'
$values[7] = 'Amén'
$values[8] = 'Amen amen gloria a Dios'
$values[9] = 'INFINITAS GRACIAS SEÑOR JESÚS por un día más de vida Por todas las maravillas y bendiciones  que nos das a diario'
$values[10] = 'Move in us Holy Spirit, fill our soul with your beautiful presence.'
$values[11] = 'Thank you my eternal GOD because your favor reaches us, every day and your great mercy there is no way to repay what you do for us, why I recognize that day by d… See more'
$values[12] = 'Amén amén aleluya yo sin ti no soy nada padre celestial ayudame señor de señores eres mi rey mi salvador te amo tanto que quiero estar junto a ti padre celestial amén '
$values[13] = 'You Are Not Alone@'
$values[14] = 'Hello blessings good day I invite you to follow this page as it has been of great blessing to my life'
$values[15] = 'Gracias gracias gracias por llegar bien a casa y librarme de todo mal amén gloria a Dios'
$values[16] = 'Amén, amén y amén'
$values[17] = 'Amén amado Dios gracias por tus misericordias.'
$values[18] = 'Amén Padre Bendito '
$values[19] = 'Amén amén gracias padre amado Dios'
$values[20] = 'Igualmente señor pongo todo a tus mano gloria a dios'
$values[21] = 'Amen gracias love'
$values[22] = 'Amén Padre Bendito '
$values[23] = 'Good Morning Happy Friday '
$values[24] = 'Amén  así es gracias padre celestial'
$values[25] = 'Bendecido viernes..'
$values[26] = 'Bendiciones..'
$values[27] = 'Amén'
$values[28] = 'Amen amen amen'
$values[29] = 'Amén '
$values[30] = 'Amen  Amen  Amen'
$values[31] = 'Gracias padre celestial amén'
$values[32] = 'Amén amén padres celestial x tdo tu bendición y también cuidados a mis hijos a cada uno de ellos'
$values[33] = 'Amén Padre Bendito '
$values[34] = 'Amén padre amado amén amén'
$values[35] = 'Amen Amen Dios siempre  dándome  una oportunidad  más.  Levantarnos de la cama bendito Dios  de verdad  por ese gran amor que nos  tiene'
$values[36] = 'Amén Padre Bendito '
$values[37] = 'Amen amen Amen'
$values[38] = 'Amén, Amén, Amén'
$values[39] = 'AMEN AMEN'
$values[40] = 'Gloria a Dios Amén y Amén '
$values[41] = 'Amen gloria a Dios'
$values[42] = 'Amén Gloria al Padre Celestial x toda su misericordia y protección .'
$values[43] = 'Amen 
 si mi señor danos fuerzas Para todos ayudanos en estos momentos mas aorita k van aoperar ami mama k ya es mayor de edad Para k salga bien de su operacion telo pedimos señor en El nombre de nuetro seño jesu crito amen 
 '
$values[44] = 'Amen amen'
$values[45] = 'Amén y Amén '
$values[46] = 'Amén Padre Bendito '
$values[47] = 'Amén amén amén'
$values[48] = 'Amén padre'
$values[49] = 'solo tu nos da un amor incondicional gracias padre  Celestial gracias por tu inmenso Amor Aleluya Gloria Dios'
$values[50] = 'Amén gloria a Dios bendiciones'
$values[51] = 'Amén amén amén'
$values[52] = 'Gracias gracias dios amén amén'
$values[53] = 'AMÉN AMÉN Y AMÉN GRACIAS GRACIAS Y GRACIAS……'
$values[54] = 'Amén aleluya gloria a Dios'
$values[55] = 'AMÉN AMÉN'
$values[56] = 'AMÉN  GRACIAS INFINITAS SEÑOR JESÚS  '
$values[57] = 'amen Gloria adios amen'
$values[58] = 'Gloria a Dios Amen padre Bendito'
$values[59] = 'Amén padre amado'
$values[60] = 'Amen amen'
$values[61] = 'Así es padre celestial amen'
$values[62] = 'Ame amen gloria a Dios'
$values[63] = 'Amén amén aleluya aleluya gloria a Dios'
$values[64] = 'Amén Amén'
$values[65] = 'Amén gloria ti señor Jesucristo Amén muchas gracias AMEN'
$values[66] = 'Dios es real'
$values[67] = 'Un Joven Cristiano.'
$values[68] = 'Amen gloria a dios padre celestial'
$values[69] = 'Un Joven Cristiano.'

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}

# Remove now-unused trailing rows (old dimension was A1:A73, new is A1:A69)
for ($row = 70; $row -le 73; $row++) {
    $ws.Cells.Item($row, 1).ClearContents()
}
